$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Value = '93_referral_statement'
$ws.Range("F18").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F19").Value = 'ppe'
$ws.Range("F21").Value = 'application instructions'
$ws.Range("F22").Value = 'env warning - water'
$ws.Range("F23").Value = 'env warning - water || off target movement'
$ws.Range("F29").Value = 'off target movement'
$ws.Range("F30").Value = 'off target movement'
$ws.Range("F31").Value = 'off target movement'
$ws.Range("F32").Value = 'off target movement'
$ws.Range("F33").Value = 'off target movement'
$ws.Range("F34").Value = 'off target movement'
$ws.Range("F35").Value = 'off target movement'
$ws.Range("F36").Value = 'off target movement'
$ws.Range("F37").Value = 'off target movement'
$ws.Range("F38").Value = 'off target movement'
$ws.Range("F39").Value = 'off target movement'
$ws.Range("F40").Value = 'off target movement'
$ws.Range("F41").Value = 'off target movement'
$ws.Range("F42").Value = 'off target movement'
$ws.Range("F43").Value = 'off target movement'
$ws.Range("F44").Value = 'off target movement'
$ws.Range("F45").Value = 'env warning - species'
$ws.Range("F47").Value = '135_product_information'
$ws.Range("F48").Value = '135_product_information'
$ws.Range("F49").Value = '135_product_information'
$ws.Range("F54").Value = 'use restrictions'
$ws.Range("F55").Value = 'use restrictions'
$ws.Range("F57").Value = 'use restrictions'
$ws.Range("F58").Value = 'use restrictions'
$ws.Range("F65").Value = 'application instructions'
$ws.Range("F66").Value = 'application instructions'
$ws.Range("F67").Value = 'application instructions'
$ws.Range("F68").Value = '134_non-agriculture_use_requirements'
$ws.Range("F70").Value = 'application instructions'
$ws.Range("F71").Value = 'application instructions'
$ws.Range("F72").Value = 'application instructions'
$ws.Range("F73").Value = 'application instructions'
$ws.Range("F74").Value = 'application instructions'
$ws.Range("F75").Value = 'application instructions'
$ws.Range("F76").Value = 'application instructions'
$ws.Range("F77").Value = 'mixing'
$ws.Range("F78").Value = 'mixing'
$ws.Range("F79").Value = 'mixing'
$ws.Range("F80").Value = 'mixing'
$ws.Range("F81").Value = 'mixing'
$ws.Range("F82").Value = 'application instructions'
$ws.Range("F83").Value = 'application instructions'
$ws.Range("F84").Value = 'application instructions'
$ws.Range("F92").Value = 'mixing'
$ws.Range("F93").Value = 'mixing'
$ws.Range("F94").Value = 'mixing'
$ws.Range("F95").Value = 'mixing'
$ws.Range("F96").Value = 'mixing'
$ws.Range("F97").Value = 'mixing'
$ws.Range("F98").Value = 'mixing'
$ws.Range("F99").Value = 'mixing'
$ws.Range("F100").Value = 'mixing'
$ws.Range("F101").Value = 'mixing'
$ws.Range("F102").Value = 'mixing'
$ws.Range("F103").Value = 'mixing'
$ws.Range("F104").Value = 'mixing'
$ws.Range("F105").Value = 'mixing'
$ws.Range("F106").Value = 'mixing'
$ws.Range("F108").Value = 'application instructions'
$ws.Range("F109").Value = 'application instructions'
$ws.Range("F112").Value = 'use restrictions'
$ws.Range("F113").Value = 'use restrictions'
$ws.Range("F114").Value = 'use restrictions'
$ws.Range("F115").Value = 'use restrictions'
$ws.Range("F116").Value = 'use restrictions'
$ws.Range("F117").Value = 'use restrictions'
$ws.Range("F118").Value = 'use restrictions'
$ws.Range("F119").Value = 'use restrictions'
$ws.Range("F120").Value = 'use restrictions'
$ws.Range("F131").Value = 'application instructions'
$ws.Range("F132").Value = 'mixing'
$ws.Range("F133").Value = 'use restrictions'
$ws.Range("F136").Value = 'application instructions'
$ws.Range("F137").Value = 'application instructions'
$ws.Range("F138").Value = 'application instructions'
$ws.Range("F140").Value = 'use restrictions'
$ws.Range("F141").Value = 'application instructions'
$ws.Range("F142").Value = 'application instructions'
$ws.Range("F143").Value = 'application instructions'
$ws.Range("F144").Value = 'use restrictions'
$ws.Range("F145").Value = 'application instructions'
$ws.Range("F146").Value = 'application instructions'
$ws.Range("F147").Value = 'use restrictions'
$ws.Range("F148").Value = 'application instructions'
$ws.Range("F149").Value = 'safety procedures'
$ws.Range("F150").Value = 'safety procedures'
$ws.Range("F151").Value = 'safety procedures'
$ws.Range("F152").Value = 'safety procedures'
$ws.Range("F153").Value = 'safety procedures'
$ws.Range("F154").Value = 'application instructions'
$ws.Range("F155").Value = 'application instructions'
$ws.Range("F156").Value = 'off target movement'
$ws.Range("F158").Value = '154_pesticide_storage'
